$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 3205326.97
$ws.Range("C9").Value = 499536.99
$ws.Range("D9").Value = 3704863.96
$ws.Range("E9").Value = 13.48327483527897
$ws.Range("F9").Value = 86.51672516472104
$ws.Range("G9").Value = -51.72224072276244
$ws.Range("H9").Value = -42.11622545171546
$ws.Range("I9").Value = 31955
$ws.Range("J9").Value = 1351
$ws.Range("K9").Value = 33306
$ws.Range("L9").Value = 22980
$ws.Range("M9").Value = 161.2212341166232
$ws.Range("N9").Value = 10.06902024189023

$wb.Save()
